# Applies the "evaluate problem 3 solution started" edit:
#  - Wraps several already-correct words in <w:proofErr> gramStart/gramEnd or
#    spellStart/spellEnd markers (cosmetic-only change to the OOXML; visible
#    text is unchanged) across a handful of existing paragraphs.
#  - Appends a new "4." / "a) ... / (blank) / b) " block (with one proofErr
#    pair) at the end of the "Predicting fingers" -> "3." answer, right after
#    the two blank indented paragraphs that follow it.

$d = $word.ActiveDocument

function Set-ParagraphXml {
    param(
        [string]$anchorText,
        [string]$innerXml
    )
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "ANCHOR NOT FOUND: $anchorText"
        return
    }
    $para = $rng.Paragraphs(1).Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $para.InsertXML($xml)
}

# 1) "Alberto Perales " -> "Alberto " + proofErr(spellStart/"Perales"/spellEnd) + " "
Set-ParagraphXml "Alberto Perales" ('<w:p><w:r><w:t xml:space="preserve">Alberto </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Perales</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>')

# 2) "b) sub goals are to find right pairing of animals and seed withing pairs in boat "
Set-ParagraphXml "sub goals are to find right pairing" ('<w:p><w:r><w:t xml:space="preserve">b) </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/><w:r><w:t>sub</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> goals are to find right pairing of animals and seed </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>withing</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> pairs in boat </w:t></w:r></w:p>')

# 3) "a) find pairing to travel ….Cat and man , seed and man, parrot and man "
Set-ParagraphXml "find pairing to travel" ('<w:p><w:r><w:t xml:space="preserve">a) </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/><w:r><w:t>find</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> pairing to travel &#8230;.Cat and man , seed and man, parrot and man </w:t></w:r></w:p>')

# 4) "a) found potential solution to pair traveling and they meet goal of not being left with wrong pair "
Set-ParagraphXml "found potential solution to pair traveling" ('<w:p><w:r><w:t xml:space="preserve">a) </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/><w:r><w:t>found</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> potential solution to pair traveling and they meet goal of not being left with wrong pair </w:t></w:r></w:p>')

# 5) "b) tried to take one at a time but it would leave impossible pairing on either side. "
Set-ParagraphXml "tried to take one at a time" ('<w:p><w:r><w:t xml:space="preserve">b) </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/><w:r><w:t>tried</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> to take one at a time but it would leave impossible pairing on either side. </w:t></w:r></w:p>')

# 6) "b)" + " his solution to find the number in each sock will work for all cases and all colors."
Set-ParagraphXml "his solution to find the number in each sock" ('<w:p><w:r><w:t>b)</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/><w:r><w:t>his</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> solution to find the number in each sock will work for all cases and all colors.</w:t></w:r></w:p>')

# 7) "a) The constraints are you loose a finger on a total count because u start on one opposite finger..."
#    (this paragraph is indented - keep its <w:pPr> intact)
Set-ParagraphXml "loose a finger on a total count" ('<w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>a)</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> The constraints </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">are you loose a finger on a total count because u </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/><w:r><w:t>start</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> on one opposite finger from both ways when counting </w:t></w:r></w:p>')

# 8) New "4." / "a) each ... / (blank) / b) " block appended after the two blank
#    indented paragraphs that follow "3. a) The solution to the sub problem...".
$rng = $d.Content
$found = $rng.Find.Execute("count from one form the first finger", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    Write-Host "ANCHOR NOT FOUND for insertion point"
} else {
    $answerPara = $rng.Paragraphs(1).Range
    # the paragraph right after the answer, then the one after that: two blank
    # indented paragraphs; insert right after the second one (before the 3rd
    # blank paragraph that originally followed).
    $blank1 = $answerPara.Next(4, 1)
    $blank2 = $blank1.Next(4, 1)
    $targetPara = $blank2.Next(4, 1)
    $targetPara = $d.Paragraphs($targetPara.Paragraphs(1).Range.Start).Range
    $insertXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>4.</w:t></w:r></w:p>' + `
        '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>a)</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>each</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> solution meets the sub goal to find what finger all counts will land on </w:t></w:r></w:p>' + `
        '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/></w:pPr></w:p>' + `
        '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t xml:space="preserve">b) </w:t></w:r></w:p>'
    $blank2.Paragraphs(1).Range.InsertXML($insertXml) | Out-Null
}

Write-Host "Done"
